# 6.4.1.2.xlsx — add a "2023" data column (Q) to the table, mirroring the
# formatting of the existing "2022" column (P), and tighten up row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow columns A:C slightly (37.140625 -> 36.28515625 chars) ---------
for ($c = 1; $c -le 3; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 35.5
}

# --- New column Q: copy formatting from column P for the rows that --------
# --- actually hold data (row 3 header, rows 5-25 body; row 4 has none) ----
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("P5:P25").Copy()
$ws.Range("Q5:Q25").PasteSpecial(-4122)

# --- Header: 2023 -----------------------------------------------------
$ws.Range("Q3").Value = 2023

# --- "В процентах к общему объему забора воды" / absolute block (rows 5-14)
$ws.Range("Q5").Value = 2385.9
$ws.Range("Q6").Value = 112.1
$ws.Range("Q7").Value = 267.89999999999998
$ws.Range("Q8").Value = 230.9
$ws.Range("Q9").Value = 249.7
$ws.Range("Q10").Value = 287
$ws.Range("Q11").Value = 334.7
$ws.Range("Q12").Value = 851
$ws.Range("Q13").Value = 48.5
$ws.Range("Q14").Value = 4.2
# Q15 stays empty (section header row, like P15)

# --- Percentage block (rows 16-25) ----------------------------------------
$ws.Range("Q16").Value = 26.890545708088244
$ws.Range("Q17").Value = 15.490056759274875
$ws.Range("Q18").Value = 22.218388220841799
$ws.Range("Q19").Value = 29.614327895683314
$ws.Range("Q20").Value = 30.104452089276922
$ws.Range("Q21").Value = 21.825966598728439
$ws.Range("Q22").Value = 32.351574864874735
$ws.Range("Q23").Value = 30.810022297218843
$ws.Range("Q24").Value = 29.193884213235311
$ws.Range("Q25").Value = 7.4362892319581295

# --- Rows 4-25 get an explicit 15pt custom height -------------------------
for ($r = 4; $r -le 25; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}
